$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 52: capitalize the three magnetic/compass terms
$ws.Range("E52").Value = "Magnetic deviation"
$ws.Range("F52").Value = "Magnetic variation"
$ws.Range("G52").Value = "Compass variation"

# Row 68: replace question text about prime meridian / date line / equator
$ws.Range("C68").Value = "They are used to measure east or west angular distance from the Prime Meridian"

# Row 69: capitalize Prime Meridian and International Date Line
$ws.Range("E69").Value = "The Prime Meridian and International Date Line are on the same great circle"

# Row 85: fix typo "within flying within" -> "flying within"
$ws.Range("F85").Value = "When flying within the signal beam"

# Row 104: prepend "to " to the three purpose phrases
$ws.Range("E104").Value = "to calculate the necessary MH to account for wind"
$ws.Range("F104").Value = "to estimate time a flight will take"
$ws.Range("G104").Value = "to create a plan to aid with pilotage"
